$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column K (2023 data). Rows not listed either
# keep K empty (label/header rows) but still need the copied number format
# from column J so the border/alignment matches the rest of the table.
$kValues = @{
    4  = 2023
    5  = 43.6
    7  = 43.5
    8  = 43.6
    10 = 52.6
    11 = 38.2
    12 = 33.9
    14 = 36.8
    15 = 47.4
    17 = 69.3
    18 = 54.9
    19 = 45.4
    20 = 49.3
    21 = 37.1
    22 = 41.1
    23 = 35.7
    24 = 36.8
    25 = 26
}

# Rows where K stays empty (section header / blank rows) but should still
# pick up J's number format (border, alignment) like the rest of the column.
$kBlankRows = @(6, 9, 13, 16)

for ($r = 4; $r -le 25; $r++) {
    $jCell = $ws.Cells.Item($r, 10)
    $kCell = $ws.Cells.Item($r, 11)
    $jCell.Copy()
    $kCell.PasteSpecial(-4122)
    if ($kValues.ContainsKey($r)) {
        $kCell.Value = $kValues[$r]
    }
}

# Clear the lingering selection highlight left over from copy/paste and park
# the cursor back on A1, matching the saved view state.
$excel.CutCopyMode = $false
$ws.Range("A1").Select()
